# Update "想去人数" (attendance count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2785
    $ws.Range("F5").Value = 6673
    $ws.Range("F6").Value = 1540
    $ws.Range("F10").Value = 87
}
